# Fix incorrect test value in test_cross_sheet_join.xlsx: 'Charles' -> 'Carol'.
# Also reflects the accompanying view-state changes captured in the commit:
#  - "customers" becomes the active/selected sheet (was "products")
#  - the active selection on "customers" moves from B1 to D7

$wb = $excel.ActiveWorkbook

$customers = $wb.Worksheets.Item("customers")

# Correct the mis-typed customer name.
$customers.Range("B4").Value = "Carol"

# Make "customers" the active sheet/tab (activeTab 1 -> 0, tabSelected swap).
$customers.Activate()

# Move the active cell/selection on "customers" from B1 to D7.
$customers.Range("D7").Select() | Out-Null
